$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (row 2: "municipio" / "Nº DE CASOS" / "Óbitos confirmados"),
# shifting all the municipio data rows up by one.
$ws.Rows(2).Delete()

# Remove the trailing two rows ("outros paises" and "outros estados"), which after
# the shift above now sit at rows 40 and 41.
$ws.Rows(40).Delete()
$ws.Rows(40).Delete()
